$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (shared string) without leaving a lingering
# custom number-format style behind (matches this workbook's convention of
# storing every Bills-sheet cell as inline text, with no style index).
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Stocks sheet: quantities on hand drop after the two new bills consume
# stock (Quantity column D).
# ---------------------------------------------------------------------------
$stocks = $wb.Worksheets.Item("Stocks")
$stocks.Range("D2").Value = 116
$stocks.Range("D3").Value = 399
$stocks.Range("D6").Value = 97

# ---------------------------------------------------------------------------
# Bills sheet: two new receipts recorded on 16-Dec-2020.
# ---------------------------------------------------------------------------
$bills = $wb.Worksheets.Item("Bills")

# Row 164 - GuestCust 33
Set-TextCell $bills 164 1 "16-Dec-2020 19:25"
Set-TextCell $bills 164 2 "GuestCust 33"
Set-TextCell $bills 164 3 "9655909777"
Set-TextCell $bills 164 4 "600"
Set-TextCell $bills 164 5 "0.0"
Set-TextCell $bills 164 6 "XX1612158"
Set-TextCell $bills 164 7 "Stock Name two(3)"

# Row 165 - GuestCust 34
Set-TextCell $bills 165 1 "16-Dec-2020 19:28"
Set-TextCell $bills 165 2 "GuestCust 34"
Set-TextCell $bills 165 3 "9655909777"
Set-TextCell $bills 165 4 "150"
Set-TextCell $bills 165 5 "0.0"
Set-TextCell $bills 165 6 "XX1612159"
Set-TextCell $bills 165 7 "Newss Sss(1)"

# Row 166 - GuestCust 35
Set-TextCell $bills 166 1 "16-Dec-2020 19:28"
Set-TextCell $bills 166 2 "GuestCust 35"
Set-TextCell $bills 166 3 "9655909777"
Set-TextCell $bills 166 4 "130.0"
Set-TextCell $bills 166 5 "0.0"
Set-TextCell $bills 166 6 "XX1612160"
Set-TextCell $bills 166 7 "Stock Name One(1)"
